$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")
$ws.Select()
$ws.Columns("G:H").Insert()
$ws.Columns("G").ColumnWidth = 69.6640625
$ws.Columns("H").ColumnWidth = 28.1640625
Write-Output "done"
